$wb = $excel.ActiveWorkbook

# --- 1. Create the new "2022-Q3" sheet -------------------------------------------------
# The new sheet needs the same layout/styling as the existing "2022-Q2" sheet, so
# duplicate it. Copying right after itself leaves:
#   [..., "2022-Q2" (original, old data), "2022-Q2 (2)" (copy, old data), "2021-Q2", ...]
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($null, $q2)

# The original "2022-Q2" sheet becomes the new "2022-Q3" sheet (with fresh values below);
# the copy keeps the old data and becomes the sheet that stays named "2022-Q2".
$q2.Name = "2022-Q3"
$q2copy = $wb.Worksheets.Item($q2.Index + 1)
$q2copy.Name = "2022-Q2"

# Restore "2021-Q2" as the active sheet/tab (matches the original workbook state).
$wb.Worksheets.Item("2021-Q2").Activate()

# --- 2. Populate "2022-Q3" with the new quarter's figures -------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue $q3.Range("D2") "20.44"
Set-TextValue $q3.Range("E2") "90.19"
Set-TextValue $q3.Range("F2") "6.21"
Set-TextValue $q3.Range("G2") "1.2693"
$q3.Range("H2").Value = 5

Set-TextValue $q3.Range("D3") "15.02"
Set-TextValue $q3.Range("E3") "90.19"
Set-TextValue $q3.Range("F3") "6.21"
Set-TextValue $q3.Range("G3") "0.9327"
$q3.Range("H3").Value = 5

# --- 3. Update the "总计" (summary) sheet -----------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Range("B2").Value = "2022-Q3"
$summary.Range("D2").Value = 2.2

$summary.Range("B3").Value = "2022-Q2"
$summary.Range("D3").Value = 2.58

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2021-Q2"
$summary.Range("C4").Value = 2
$summary.Range("D4").Value = 2.94

$summary.Range("A3").Copy()
$summary.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
